{"js": "// Benchmark stats table update: the single table in the document has one\n// column; each row holds one metric value. A handful of rows get their\n// text replaced, and the final three rows (which each still carried the\n// whole tab-separated \"db-shootout\" line bundled into a single run) get\n// collapsed down to just their first field.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// index (0-based) -> new text for that row's single cell\nconst replacements = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"5212\",\n  4: \"0.00001\",\n  5: \"0.02179\",\n  6: \"0.00014\",\n  7: \"0.00016\",\n  8: \"0.00017\",\n  9: \"0.00019\",\n  10: \"0.00026\",\n  11: \"0.78900\",\n  43: \"99.89\",\n  44: \"0.79\",\n  45: \"693\",\n};\n\nfor (const [idxStr, newText] of Object.entries(replacements)) {\n  const idx = Number(idxStr);\n  const cell = rows.items[idx].cells.items[0];\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Benchmark stats table update: the single table in the document has one\n# column; each row holds one metric value. A handful of rows get their\n# text replaced, and the final three rows (which each still carried the\n# whole tab-separated \"db-shootout\" line bundled into a single run) get\n# collapsed down to just their first field.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based row index -> new cell text\n$replacements = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"5212\"\n    5  = \"0.00001\"\n    6  = \"0.02179\"\n    7  = \"0.00014\"\n    8  = \"0.00016\"\n    9  = \"0.00017\"\n    10 = \"0.00019\"\n    11 = \"0.00026\"\n    12 = \"0.78900\"\n    44 = \"99.89\"\n    45 = \"0.79\"\n    46 = \"693\"\n}\n\nforeach ($rowIndex in $replacements.Keys) {\n    $t.Cell($rowIndex, 1).Range.Text = $replacements[$rowIndex]\n}\n"}
